# Scheduled-runner update: refresh market-board derived price/profit
# columns (currentAveragePrice*, LevePrice*, LeveProfit*) across the
# per-job leve-profit sheets. Values only; no structural changes.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4642.857
$ws.Range("J64").Value = 4642.857
$ws.Range("L64").Value = 4642.857
$ws.Range("N64").Value = -5138.857
$ws.Range("H67").Value = 4642.857
$ws.Range("J67").Value = 4642.857
$ws.Range("L67").Value = 4642.857
$ws.Range("N67").Value = -6358.857
$ws.Range("H129").Value = 264329.72
$ws.Range("I129").Value = 287.5
$ws.Range("J129").Value = 295393.5
$ws.Range("K129").Value = 862.5
$ws.Range("L129").Value = 886180.5
$ws.Range("M129").Value = 4137.5
$ws.Range("N129").Value = -896180.5
$ws.Range("H135").Value = 11908176
$ws.Range("I135").Value = 708.06665
$ws.Range("K135").Value = 6372.59985
$ws.Range("M135").Value = -3837.59985
$ws.Range("H137").Value = 1766.8286
$ws.Range("I137").Value = 1608.2413
$ws.Range("K137").Value = 4824.7239
$ws.Range("M137").Value = -2274.7239
$ws.Range("H138").Value = 10991196
$ws.Range("I138").Value = 21739972
$ws.Range("J138").Value = 3557.3333
$ws.Range("K138").Value = 65219916
$ws.Range("L138").Value = 10671.9999
$ws.Range("M138").Value = -65214776
$ws.Range("N138").Value = -20951.9999
$ws.Range("H141").Value = 1144.4667
$ws.Range("I141").Value = 858.46344
$ws.Range("J141").Value = 4076
$ws.Range("K141").Value = 2575.39032
$ws.Range("L141").Value = 12228
$ws.Range("M141").Value = 2604.60968
$ws.Range("N141").Value = -22588

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2012.15
$ws.Range("I32").Value = 2047.4286
$ws.Range("J32").Value = 1655.4445
$ws.Range("K32").Value = 2047.4286
$ws.Range("L32").Value = 1655.4445
$ws.Range("M32").Value = -1760.4286
$ws.Range("N32").Value = -2229.4445
$ws.Range("H63").Value = 2102.5
$ws.Range("I63").Value = 2102.5
$ws.Range("K63").Value = 2102.5
$ws.Range("M63").Value = -1416.5
$ws.Range("H66").Value = 2102.5
$ws.Range("I66").Value = 2102.5
$ws.Range("K66").Value = 10512.5
$ws.Range("M66").Value = -7080.5
$ws.Range("H74").Value = 37039644
$ws.Range("I74").Value = 43480810
$ws.Range("J74").Value = 2949.5
$ws.Range("K74").Value = 43480810
$ws.Range("L74").Value = 2949.5
$ws.Range("M74").Value = -43479936
$ws.Range("N74").Value = -4697.5
$ws.Range("H77").Value = 37039644
$ws.Range("I77").Value = 43480810
$ws.Range("J77").Value = 2949.5
$ws.Range("K77").Value = 217404050
$ws.Range("L77").Value = 14747.5
$ws.Range("M77").Value = -217399682
$ws.Range("N77").Value = -23483.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3199.6428
$ws.Range("J105").Value = 2650
$ws.Range("L105").Value = 2650
$ws.Range("N105").Value = -6144
$ws.Range("H141").Value = 54780
$ws.Range("J141").Value = 54780
$ws.Range("L141").Value = 54780
$ws.Range("N141").Value = -65140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H31").Value = 3344.0264
$ws.Range("I31").Value = 1931.8846
$ws.Range("K31").Value = 1931.8846
$ws.Range("M31").Value = -1636.8846
$ws.Range("H34").Value = 3344.0264
$ws.Range("I34").Value = 1931.8846
$ws.Range("K34").Value = 1931.8846
$ws.Range("M34").Value = -1729.8846
$ws.Range("H62").Value = 33336476
$ws.Range("I62").Value = 38464336
$ws.Range("J62").Value = 5376.5
$ws.Range("K62").Value = 38464336
$ws.Range("L62").Value = 5376.5
$ws.Range("M62").Value = -38463712
$ws.Range("N62").Value = -6624.5
$ws.Range("H65").Value = 33336476
$ws.Range("I65").Value = 38464336
$ws.Range("J65").Value = 5376.5
$ws.Range("K65").Value = 192321680
$ws.Range("L65").Value = 26882.5
$ws.Range("M65").Value = -192318560
$ws.Range("N65").Value = -33122.5
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("M124").ClearContents()
$ws.Range("N124").ClearContents()
$ws.Range("H125").Value = 19080
$ws.Range("J125").Value = 27966.666
$ws.Range("L125").Value = 27966.666
$ws.Range("N125").Value = -32886.666
$ws.Range("H132").Value = 2055.2
$ws.Range("I132").Value = 1589.3158
$ws.Range("K132").Value = 4767.9474
$ws.Range("M132").Value = -2237.9474

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 3749.75
$ws.Range("J107").Value = 188.4375
$ws.Range("L107").Value = 565.3125
$ws.Range("N107").Value = -4405.3125
$ws.Range("H113").Value = 682.3200000000001
$ws.Range("I113").Value = 513.3333
$ws.Range("J113").Value = 935.8
$ws.Range("K113").Value = 1539.9999
$ws.Range("L113").Value = 2807.4
$ws.Range("M113").Value = 630.0001
$ws.Range("N113").Value = -7147.4
$ws.Range("H116").Value = 1685.7142
$ws.Range("I116").Value = 1500
$ws.Range("J116").Value = 1933.3334
$ws.Range("K116").Value = 4500
$ws.Range("L116").Value = 5800.0002
$ws.Range("M116").Value = -1058
$ws.Range("N116").Value = -12684.0002
$ws.Range("H131").Value = 737.6429000000001
$ws.Range("J131").Value = 748.37634
$ws.Range("L131").Value = 2245.12902
$ws.Range("N131").Value = -12325.12902
$ws.Range("H132").Value = 574.5
$ws.Range("J132").Value = 999
$ws.Range("L132").Value = 8991
$ws.Range("N132").Value = -14051
$ws.Range("H133").Value = 3460.2222
$ws.Range("I133").Value = 2676.6667
$ws.Range("J133").Value = 3852
$ws.Range("K133").Value = 8030.000100000001
$ws.Range("L133").Value = 11556
$ws.Range("M133").Value = -2970.000100000001
$ws.Range("N133").Value = -21676
$ws.Range("H139").Value = 1529.1936
$ws.Range("I139").Value = 1042.8846
$ws.Range("J139").Value = 4058
$ws.Range("K139").Value = 3128.6538
$ws.Range("L139").Value = 12174
$ws.Range("M139").Value = 2011.3462
$ws.Range("N139").Value = -22454

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 178222.22
$ws.Range("J24").Value = 102000
$ws.Range("L24").Value = 102000
$ws.Range("N24").Value = -102346

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 3335.6667
$ws.Range("J24").Value = 3335.6667
$ws.Range("L24").Value = 3335.6667
$ws.Range("N24").Value = -4021.6667
$ws.Range("H56").Value = 16000
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 16000
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 16000
$ws.Range("M56").ClearContents()
$ws.Range("N56").Value = -17382
$ws.Range("H132").Value = 1680.1
$ws.Range("I132").Value = 1680.1
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5040.299999999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2510.299999999999
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 1582.125
$ws.Range("I136").Value = 1464.5714
$ws.Range("J136").Value = 2405
$ws.Range("K136").Value = 4393.7142
$ws.Range("L136").Value = 7215
$ws.Range("M136").Value = -1843.7142
$ws.Range("N136").Value = -12315

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 14683.667
$ws.Range("J61").Value = 18000
$ws.Range("L61").Value = 18000
$ws.Range("N61").Value = -18584
$ws.Range("H62").Value = 4600.4
$ws.Range("I62").Value = 4001
$ws.Range("K62").Value = 4001
$ws.Range("M62").Value = -3377
$ws.Range("H65").Value = 4600.4
$ws.Range("I65").Value = 4001
$ws.Range("K65").Value = 20005
$ws.Range("M65").Value = -16885
$ws.Range("H132").Value = 803
$ws.Range("I132").Value = 803
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2409
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = 121
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 17859310
$ws.Range("I136").Value = 25642050
$ws.Range("K136").Value = 76926150
$ws.Range("M136").Value = -76923600
